$wb = $excel.ActiveWorkbook

# Rename sheets (append "_1" suffix)
$wb.Worksheets.Item("PR_DEN_1_Pipettes").Name = "PR_DEN_1_Pipettes_1"
$wb.Worksheets.Item("PR_DEN_2_Pipettes").Name = "PR_DEN_2_Pipettes_1"
$wb.Worksheets.Item("PR_DEN_1_Cassette").Name = "PR_DEN_1_Cassette_1"
$wb.Worksheets.Item("PR_DEN_3_2_Cassette").Name = "PR_DEN_3_2_Cassette_1"

# The rename does not rewrite the stored Print_Area defined-name formula text,
# so re-apply each renamed sheet's print area to refresh the sheet-qualified reference.
$wb.Worksheets.Item("PR_DEN_1_Pipettes_1").PageSetup.PrintArea = "A1:G42"
$wb.Worksheets.Item("PR_DEN_2_Pipettes_1").PageSetup.PrintArea = "A1:G42"
$wb.Worksheets.Item("PR_DEN_1_Cassette_1").PageSetup.PrintArea = "A1:H42"
$wb.Worksheets.Item("PR_DEN_3_2_Cassette_1").PageSetup.PrintArea = "A1:H42"

# Adjust row 11 height on the sheets that still need it
$wb.Worksheets.Item("PR_DEN_2_Uncut_Sheet_1").Rows.Item(11).RowHeight = 15.25
$wb.Worksheets.Item("PR_DEN_2_Pipettes_1").Rows.Item(11).RowHeight = 15.25
$wb.Worksheets.Item("PR_DEN_3_2_Cassette_1").Rows.Item(11).RowHeight = 15.25

# Move the active tab from PR_DEN_Buffer_2 to PR_DEN_3_2_Cassette_1, and change its selection
$wb.Worksheets.Item("PR_DEN_3_2_Cassette_1").Select()
$wb.Worksheets.Item("PR_DEN_3_2_Cassette_1").Range("B6").Select()
